# Apply the "nuevos experimentos no convexos" numeric updates.
#
# All of the edited cells in this workbook hold their (numeric-looking)
# values as plain text (shared strings), so a plain `.Value = "..."`
# assignment must be prevented from being auto-coerced into a real
# number by Excel. We force text by temporarily setting the cell's
# NumberFormat to "@" (Text) before assigning the value, then clear the
# format again so the cell keeps the workbook's default (General) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# Sheet order (per workbook.xml):
#   1 Funciones_Objetivo            (unaffected)
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha                  (unaffected)
#
# NOTE: sheets 5 and 6 are named "Vector_bf" / "Vector_BF" which differ
# only by case, so Worksheets.Item(name) would be ambiguous. Using the
# numeric (1-based) index avoids that problem entirely.

$wsLider = $wb.Worksheets.Item(2)
Set-TextValue $wsLider.Range("A2") "-0.9 + x"
Set-TextValue $wsLider.Range("B2") "-0.09999999999999998"
Set-TextValue $wsLider.Range("D2") "0.42"
Set-TextValue $wsLider.Range("A3") "0.8999999999999999 - x"
Set-TextValue $wsLider.Range("B3") "-1.9"
Set-TextValue $wsLider.Range("D3") "0.02"

$wsFollower = $wb.Worksheets.Item(3)
Set-TextValue $wsFollower.Range("A2") "-2.7 + y"
Set-TextValue $wsFollower.Range("B2") "1.7000000000000002"
Set-TextValue $wsFollower.Range("D2") "0.29"
Set-TextValue $wsFollower.Range("E2") "7.5"
Set-TextValue $wsFollower.Range("F2") "9.200000000000001"
Set-TextValue $wsFollower.Range("A3") "2.7 - y"
Set-TextValue $wsFollower.Range("B3") "-3.7"
Set-TextValue $wsFollower.Range("D3") "0.52"
Set-TextValue $wsFollower.Range("E3") "6.8999999999999995"
Set-TextValue $wsFollower.Range("F3") "7.4"

$wsPunto = $wb.Worksheets.Item(4)
Set-TextValue $wsPunto.Range("A2") "0.9"
Set-TextValue $wsPunto.Range("B2") "2.7"

$wsVecbf = $wb.Worksheets.Item(5)
Set-TextValue $wsVecbf.Range("A2") "-6.16"

$wsVecBF = $wb.Worksheets.Item(6)
Set-TextValue $wsVecBF.Range("A2") "-1.7"
Set-TextValue $wsVecBF.Range("A3") "-6.000000000000001"
